$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column (H1), reusing the header formatting from the
# existing "sum" header cell (G1) so no duplicate style is introduced.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for the data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
